$d = $word.ActiveDocument

$replacements = @(
    @("2025-04-01 Tuesday", "2025-04-02 Wednesday"),
    @("236÷6=", "995÷2="),
    @("649÷9=", "354÷2="),
    @("709÷4=", "271÷2="),
    @("876÷2=", "956÷8="),
    @("793÷4=", "371÷4="),
    @("303÷9=", "671÷6="),
    @("650÷8=", "901÷6="),
    @("239÷7=", "496÷8="),
    @("457÷4=", "606÷6="),
    @("757÷5=", "752÷7="),
    @("218÷6=", "796÷2="),
    @("751÷6=", "753÷5="),
    @("371÷2=", "896÷7="),
    @("848÷4=", "691÷4="),
    @("859÷6=", "123÷5="),
    @("904÷4=", "699÷3="),
    @("366÷6=", "770÷6="),
    @("540÷9=", "615÷2="),
    @("551÷5=", "331÷8="),
    @("383÷9=", "408÷4="),
    @("129÷9=", "169÷8="),
    @("280÷7=", "373÷4="),
    @("995÷7=", "287÷5="),
    @("779÷5=", "148÷6="),
    @("698÷2=", "392÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
